$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$headerRange = $ws.Range("A1:U1")
$headerRange.ClearFormats()
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U88"), $null, 1)
$headerRange.Font.Bold = $true
$headerRange.Interior.Pattern = 1
$headerRange.Interior.Color = 14277081
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.WrapText = $true
Write-Output "done"
